$win = $excel.ActiveWindow
Write-Host ("Zoom=" + $win.Zoom)
$win.WindowState = -4137  # xlMaximized
